$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TP")
$ws.Range("B2").Value = "test"
